$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Append the new log row (row 4) mirroring the structure of rows 2/3.
$logs.Range("A4").Value = "Retour status"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("D4").Value = "Retour / Terugbetaling"
$logs.Range("F4").Value = "2025-08-28 17:50:25"
$logs.Range("G4").Value = "Nee"
$logs.Range("H4").Value = "Ja"
$logs.Range("I4").Value = "Nee"
$logs.Range("J4").Value = "Nee"

# Bump the Dashboard counter for this category.
$dash.Range("B2").Value = 3

# Extend the conditional-formatting ranges (D/G/H/I/J) to cover the new row.
$logs.Range("D2:D3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D4"))
$logs.Range("G2:G3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G4"))
$logs.Range("H2:H3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H4"))
$logs.Range("I2:I3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I4"))
$logs.Range("J2:J3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J4"))
